$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 01:08"

# Row 4: numeric updates: B4=5564920, C4=35131, D4=2918742, E4=2473106, G4=466, H4=173072
$ws.Range("B4").Value = 5564920
$ws.Range("C4").Value = 35131
$ws.Range("D4").Value = 2918742
$ws.Range("E4").Value = 2473106
$ws.Range("G4").Value = 466
$ws.Range("H4").Value = 173072

# Row 11: numeric updates: B11=468332, C11=11643, D11=287436, E11=165799, G11=287, H11=15097
$ws.Range("B11").Value = 468332
$ws.Range("C11").Value = 11643
$ws.Range("D11").Value = 287436
$ws.Range("E11").Value = 165799
$ws.Range("G11").Value = 287
$ws.Range("H11").Value = 15097

# Row 17: numeric updates: B17=294569, C17=5469, E17=77164, G17=66, H17=5703
$ws.Range("B17").Value = 294569
$ws.Range("C17").Value = 5469
$ws.Range("E17").Value = 77164
$ws.Range("G17").Value = 66
$ws.Range("H17").Value = 5703

# Row 32: numeric updates: B32=96475, C32=139, D32=59743, E32=31572, G32=19, H32=5160
$ws.Range("B32").Value = 96475
$ws.Range("C32").Value = 139
$ws.Range("D32").Value = 59743
$ws.Range("E32").Value = 31572
$ws.Range("G32").Value = 19
$ws.Range("H32").Value = 5160

# Row 35: numeric updates: B35=86309, C35=764, D35=52905, E35=31951, G35=15, H35=1453
$ws.Range("B35").Value = 86309
$ws.Range("C35").Value = 764
$ws.Range("D35").Value = 52905
$ws.Range("E35").Value = 31951
$ws.Range("G35").Value = 15
$ws.Range("H35").Value = 1453

# Row 49: country -> Japon; numeric updates: B49=54714, C49=1137, D49=40080, E49=13546, H49=1088
$ws.Range("A49").Value = "Japon"
$ws.Range("B49").Value = 54714
$ws.Range("C49").Value = 1137
$ws.Range("D49").Value = 40080
$ws.Range("E49").Value = 13546
$ws.Range("H49").Value = 1088

# Row 50: country -> Portugal; numeric updates: B50=54102, C50=121, D50=39697, E50=12627, G50=3, H50=1778
$ws.Range("A50").Value = "Portugal"
$ws.Range("B50").Value = 54102
$ws.Range("C50").Value = 121
$ws.Range("D50").Value = 39697
$ws.Range("E50").Value = 12627
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 1778

# Row 52: numeric updates: B52=49068, C52=298, D52=36497, E52=11596, G52=1, H52=975
$ws.Range("B52").Value = 49068
$ws.Range("C52").Value = 298
$ws.Range("D52").Value = 36497
$ws.Range("E52").Value = 11596
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 975

# Row 63: numeric updates: B63=33755, C63=1148, D63=22700, E63=10774, G63=5, H63=281
$ws.Range("B63").Value = 33755
$ws.Range("C63").Value = 1148
$ws.Range("D63").Value = 22700
$ws.Range("E63").Value = 10774
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 281

# Row 74: numeric updates: B74=20012, C74=121, D74=13799, E74=5816
$ws.Range("B74").Value = 20012
$ws.Range("C74").Value = 121
$ws.Range("D74").Value = 13799
$ws.Range("E74").Value = 5816

# Row 76: numeric updates: B76=17026, C76=33, E76=2969
$ws.Range("B76").Value = 17026
$ws.Range("C76").Value = 33
$ws.Range("E76").Value = 2969

# Row 81: numeric updates: B81=14365, C81=32, D81=9186, E81=4681, G81=3, H81=498
$ws.Range("B81").Value = 14365
$ws.Range("C81").Value = 32
$ws.Range("D81").Value = 9186
$ws.Range("E81").Value = 4681
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 498

# Row 86: numeric updates: B86=10005, C86=40, E86=887
$ws.Range("B86").Value = 10005
$ws.Range("C86").Value = 40
$ws.Range("E86").Value = 887

# Row 96: numeric updates: B96=7879, C96=48, E96=2448
$ws.Range("B96").Value = 7879
$ws.Range("C96").Value = 48
$ws.Range("E96").Value = 2448

# Row 107: numeric updates: B107=5072, C107=46, D107=2626, E107=2285, G107=4, H107=161
$ws.Range("B107").Value = 5072
$ws.Range("C107").Value = 46
$ws.Range("D107").Value = 2626
$ws.Range("E107").Value = 2285
$ws.Range("G107").Value = 4
$ws.Range("H107").Value = 161

# Row 114: numeric updates: B114=4035, C114=75, D114=2910, E114=1048, G114=2, H114=77
$ws.Range("B114").Value = 4035
$ws.Range("C114").Value = 75
$ws.Range("D114").Value = 2910
$ws.Range("E114").Value = 1048
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 77

# Row 143: numeric updates: B143=1440, C143=6, D143=1200
$ws.Range("B143").Value = 1440
$ws.Range("C143").Value = 6
$ws.Range("D143").Value = 1200

# Row 164: country -> Guyana; numeric updates: B164=709, C164=35, D164=349, E164=337, G164=1, H164=23
$ws.Range("A164").Value = "Guyana"
$ws.Range("B164").Value = 709
$ws.Range("C164").Value = 35
$ws.Range("D164").Value = 349
$ws.Range("E164").Value = 337
$ws.Range("G164").Value = 1
$ws.Range("H164").Value = 23

# Row 165: country -> San Marino; numeric updates: B165=699, D165=657, E165=0, H165=42
$ws.Range("A165").Value = "San Marino"
$ws.Range("B165").Value = 699
$ws.Range("D165").Value = 657
$ws.Range("E165").Value = 0
$ws.Range("H165").Value = 42

# Row 166: numeric updates: B166=552, C166=55, E166=401
$ws.Range("B166").Value = 552
$ws.Range("C166").Value = 55
$ws.Range("E166").Value = 401

# Row 213: country -> Islas Malvinas; numeric updates: D213=13, H213=0
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214: country -> Montserrat; numeric updates: D214=12, H214=1
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
